# Refresh the crypto price (D) and 1h-volume-change (E) columns.
# D holds price text that can look numeric (thousand-dot grouping like
# '26.284.24', trailing zeros like '64.50', subscript-zero notation like
# '0.0₃0727'), so the whole column is switched to text format first to
# stop Excel's auto-detection from turning these into Numbers / dropping
# significant trailing zeros, then the style is restored to Normal so the
# cells end up with no explicit style, matching the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = '26.284.24'
$ws.Range("E2").Value = '  +0.38%  '
$ws.Range("D3").Value = '1.591.00'
$ws.Range("E3").Value = '  +0.66%  '
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").Value = '213.08'
$ws.Range("E5").Value = '  +1.68%  '
$ws.Range("E7").Value = '  -0.27%  '
$ws.Range("E8").Value = '  +0.33%  '
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("D10").Value = '19.37'
$ws.Range("E10").Value = '  -0.60%  '
$ws.Range("D11").Value = '0.0849'
$ws.Range("E11").Value = '  +0.44%  '
$ws.Range("D12").Value = '1.814.85'
$ws.Range("E12").Value = '  +0.64%  '
$ws.Range("D13").Value = '1.585.60'
$ws.Range("E13").Value = '  -0.19%  '
$ws.Range("E14").Value = '  +0.10%  '
$ws.Range("E15").Value = '  +1.45%  '
$ws.Range("D16").Value = '64.50'
$ws.Range("E16").Value = '  +0.24%  '
$ws.Range("D17").Value = '26.281.66'
$ws.Range("E17").Value = '  +0.33%  '
$ws.Range("D18").Value = '0.0₃0727'
$ws.Range("E18").Value = '  -0.97%  '
$ws.Range("E19").Value = '  +2.74%  '
$ws.Range("D20").Value = '213.34'
$ws.Range("E20").Value = '  +3.01%  '
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("E22").Value = '  +0.97%  '
$ws.Range("D23").Value = '9.00'
$ws.Range("E23").Value = '  +1.40%  '
$ws.Range("D24").Value = '2.16'
$ws.Range("E24").Value = '  -2.05%  '
$ws.Range("D25").Value = '144.90'
$ws.Range("E25").Value = '  +0.31%  '
$ws.Range("E26").Value = '  -0.24%  '
$ws.Range("E27").Value = '  +0.96%  '
$ws.Range("E28").Value = '  -0.43%  '
$ws.Range("D29").Value = '15.21'
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  -0.94%  '
$ws.Range("E31").Value = '  +1.27%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D34").Value = '1.338.48'
$ws.Range("E34").Value = '  +4.86%  '
$ws.Range("E35").Value = '  -0.89%  '
$ws.Range("E36").Value = '  -0.20%  '
$ws.Range("D37").Value = '0.593'
$ws.Range("E37").Value = '  -2.83%  '
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("D39").Value = '0.820'
$ws.Range("E39").Value = '  +0.35%  '
$ws.Range("E40").Value = '  +4.29%  '
$ws.Range("E42").Value = '  -1.04%  '
$ws.Range("E43").Value = '  +0.28%  '
$ws.Range("D44").Value = '0.763'
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("D45").Value = '61.90'
$ws.Range("E45").Value = '  -0.48%  '
$ws.Range("D46").Value = '1.725.74'
$ws.Range("E46").Value = '  +0.50%  '
$ws.Range("D47").Value = '86.52'
$ws.Range("E47").Value = '  -2.81%  '
$ws.Range("D48").Value = '0.0₆0103'
$ws.Range("E48").Value = '  -1.77%  '
$ws.Range("E49").Value = '  -3.43%  '
$ws.Range("D50").Value = '0.0982'
$ws.Range("E50").Value = '  -2.06%  '
$ws.Range("E51").Value = '  -0.27%  '

$priceCol.Style = "Normal"
